$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.888.25'
$ws.Range("E2").Value = '  +1.97%  '

$ws.Range("D3").Value = '1.711.63'
$ws.Range("E3").Value = '  +2.01%  '

$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  -0.25%  '

$ws.Range("D5").Value = '313.39'
$ws.Range("E5").Value = '  +2.21%  '

$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("D7").Value = '0.3761'
$ws.Range("E7").Value = '  +1.63%  '

$ws.Range("D8").Value = '49.73'
$ws.Range("E8").Value = '  +3.04%  '

$ws.Range("D9").Value = '0.3472'
$ws.Range("E9").Value = '  +0.75%  '

$ws.Range("D10").Value = '1.222'
$ws.Range("E10").Value = '  +3.71%  '

$ws.Range("D11").Value = '0.07601'
$ws.Range("E11").Value = '  +4.67%  '

$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.07%  '

$ws.Range("E13").Value = '  +5.37%  '

$ws.Range("D14").Value = '6.363'
$ws.Range("E14").Value = '  +3.70%  '

$ws.Range("D15").Value = '7.103'
$ws.Range("E15").Value = '  +5.35%  '

$ws.Range("D16").Value = '1.710.79'
$ws.Range("E16").Value = '  +1.95%  '

$ws.Range("D17").Value = '0.00001137'
$ws.Range("E17").Value = '  +2.71%  '

$ws.Range("E18").Value = '  +0.47%  '

$ws.Range("D19").Value = '0.9990'
$ws.Range("E19").Value = '  -0.02%  '

$ws.Range("D20").Value = '85.23'
$ws.Range("E20").Value = '  +5.13%  '

$ws.Range("E21").Value = '  +6.14%  '

$ws.Range("D22").Value = '6.425'
$ws.Range("E22").Value = '  +5.37%  '

$ws.Range("D23").Value = '13.27'
$ws.Range("E23").Value = '  +10.49%  '

$ws.Range("D24").Value = '24.886.85'
$ws.Range("E24").Value = '  +2.08%  '

$ws.Range("D25").Value = '2.466'
$ws.Range("E25").Value = '  +1.19%  '

$ws.Range("D26").Value = '2.811'
$ws.Range("E26").Value = '  +5.10%  '

$ws.Range("D27").Value = '20.58'
$ws.Range("E27").Value = '  +5.67%  '

$ws.Range("D28").Value = '151.12'
$ws.Range("E28").Value = '  -0.81%  '

$ws.Range("D29").Value = '133.29'
$ws.Range("E29").Value = '  +5.20%  '

$ws.Range("D30").Value = '1.903.65'
$ws.Range("E30").Value = '  +2.24%  '

$ws.Range("D31").Value = '1.255'
$ws.Range("E31").Value = '  +28.10%  '

$ws.Range("D32").Value = '6.937'
$ws.Range("E32").Value = '  +9.15%  '

$ws.Range("D33").Value = '4.243'
$ws.Range("E33").Value = '  +4.95%  '

$ws.Range("D34").Value = '13.97'
$ws.Range("E34").Value = '  +11.07%  '

$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = '0.08890'
$ws.Range("E35").Value = '  +5.41%  '

$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '1.774'
$ws.Range("E36").Value = '  +4.35%  '

$ws.Range("D37").Value = '5.689'
$ws.Range("E37").Value = '  +6.21%  '

$ws.Range("D38").Value = '9.387'
$ws.Range("E38").Value = '  +5.39%  '

$ws.Range("D39").Value = '0.06702'
$ws.Range("E39").Value = '  +2.79%  '

$ws.Range("D40").Value = '0.02426'
$ws.Range("E40").Value = '  +4.14%  '

$ws.Range("D41").Value = '0.2248'
$ws.Range("E41").Value = '  +6.45%  '

$ws.Range("E42").Value = '  +1.98%  '

$ws.Range("D43").Value = '0.6493'
$ws.Range("E43").Value = '  +5.10%  '

$ws.Range("D44").Value = '0.9994'
$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("D45").Value = '14.01'
$ws.Range("E45").Value = '  +5.87%  '

$ws.Range("E46").Value = '  +4.23%  '

$ws.Range("D47").Value = '3.842'
$ws.Range("E47").Value = '  +2.20%  '

$ws.Range("D48").Value = '2.147'
$ws.Range("E48").Value = '  +5.88%  '

$ws.Range("D49").Value = '131.01'
$ws.Range("E49").Value = '  +2.86%  '

$ws.Range("D50").Value = '0.07333'
$ws.Range("E50").Value = '  +1.83%  '

$ws.Range("D51").Value = '80.63'
$ws.Range("E51").Value = '  +6.32%  '
